$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: email used for login test data was updated to a new address.
$ws.Range("D2").Value = "lakherar171@gmail.com"

# F2: add a hyperlink (mailto) for the new email/password pairing and
# update the password value, matching the Hyperlink cell style used
# elsewhere on the row (reuse existing style via a format-only paste
# from D2, since Hyperlinks.Add on its own introduces a brand-new xf).
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:lakherar171@gmail.com")
$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Lakhera@1998"

# Clear the clipboard/marching-ants state left behind by Copy().
$excel.CutCopyMode = 0

# The active selection moved to B1 before the workbook was saved.
[void]$ws.Range("B1").Select()
